$p = $ppt.ActivePresentation

# --- Slide 11: "User story 2 - TOTAL VOLUME TRADED with entity" ---
# Content Placeholder 2 (p:ph idx="1") currently inherits its position from the
# slide layout (empty <p:spPr/>). Bake in an explicit transform matching the
# inherited geometry, and add an extra trailing space after "Filtering : ".
$s11 = $p.Slides.Item(11)
$sh11 = $s11.Shapes.Item(2)

$sh11.Left = 114.29755905511811
$sh11.Top = 158.71905511811025
$sh11.Width = 756.1633858267717
$sh11.Height = 271.70181102362204

$tr11 = $sh11.TextFrame.TextRange
$full11 = $tr11.Text
$idx11 = $full11.IndexOf("Filtering : ")
$sub11 = $tr11.Characters($idx11 + 1, 12)
$sub11.Text = "Filtering :  "

# --- Slide 12: "User story 3 - RFQ strike rates" ---
# Same extra trailing space after "Filtering : ".
$s12 = $p.Slides.Item(12)
$sh12 = $s12.Shapes.Item(2)

$tr12 = $sh12.TextFrame.TextRange
$full12 = $tr12.Text
$idx12 = $full12.IndexOf("Filtering : ")
$sub12 = $tr12.Characters($idx12 + 1, 12)
$sub12.Text = "Filtering :  "

# --- Slide 8: "Sprint 1I: Details" title ---
$s8 = $p.Slides.Item(8)
$sh8 = $s8.Shapes.Item(1)
$sh8.TextFrame.TextRange.Text = "Sprint 1I : Details"
